$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the top of this subset's data
# block (row 191), pushing all the existing rows (191-220) down by one
# (to 192-221).
$ws.Rows.Item(191).Insert()

$ws.Cells.Item(191, 1).Value = 7
$ws.Cells.Item(191, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(191, 3).Value = "Ñuble"
$ws.Cells.Item(191, 4).Value = 44776
$ws.Cells.Item(191, 5).Value = 16
$ws.Cells.Item(191, 6).Value = 100112017
$ws.Cells.Item(191, 7).Value = "Apio"
$ws.Cells.Item(191, 8).Value = "Americana (o)"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 120
$ws.Cells.Item(191, 11).Value = 9000
$ws.Cells.Item(191, 12).Value = 10000
$ws.Cells.Item(191, 13).Value = 9500
$ws.Cells.Item(191, 14).Value = "$/docena de matas"
$ws.Cells.Item(191, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(191, 16).Value = 1583
$ws.Cells.Item(191, 17).Value = 6
$ws.Cells.Item(191, 18).Value = "Hortaliza"
